# "add user list to project"
#
# 1. On "PI hours", the cfop list for Romit Roy Choudhury (row 2, col G) is
#    reordered from ['cfop_CHOUDHURY', 'cfop_RRC'] to ['cfop_RRC', 'cfop_CHOUDHURY'].
# 2. On "project hours", a new "users" column (E) is added with a header and
#    one user list per project row.

$wb = $excel.ActiveWorkbook

# --- 1. Fix the cfop ordering for Romit Roy Choudhury on "PI hours" ---
$wsPI = $wb.Worksheets.Item("PI hours")
$wsPI.Cells.Item(2, 7).Value = "['cfop_RRC', 'cfop_CHOUDHURY']"

# --- 2. Add the "users" column to "project hours" ---
$wsProj = $wb.Worksheets.Item("project hours")

# Clone formatting (bold/border/centered header style, borders on column A
# equivalents) from column D into the new column E first, so the new cells
# pick up the same look as the rest of the header/body - then overwrite the
# values explicitly.
$wsProj.Range("D1:D5").Copy($wsProj.Range("E1:E5"))

$wsProj.Cells.Item(1, 5).Value = "users"
$wsProj.Cells.Item(2, 5).Value = "['Mahanth Gowda', 'Ashutosh Dhekne, Mahanth Gowda, Sheng Shen', 'Ashutosh Dhekne']"
$wsProj.Cells.Item(3, 5).Value = "['HYUNG JIN YOON', 'Hyung-Jin Yoon']"
$wsProj.Cells.Item(4, 5).Value = "['Joseph Chapman']"
$wsProj.Cells.Item(5, 5).Value = "['Shuchen Song']"
